$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quant")
$ws.Activate()

# Delete column A entirely (shifts B,C,D left to A,B,C)
$ws.Columns.Item(1).Delete() | Out-Null

# Set the selection to B11 as the new active cell
$ws.Range("B11").Select() | Out-Null
